# Testing redownload & backup
# The "AI" sheet's article-URL list was refreshed: the old A3:A22 block
# (20 URLs) is replaced by a new, shorter A3:A20 block (18 URLs), and the
# sheet's active selection moves back up to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AI")

$urls = @(
    "https://response.jp/article/2019/07/09/324248.html",
    "https://response.jp/article/2019/04/18/321534.html",
    "https://response.jp/article/2019/12/06/329547.html",
    "https://response.jp/article/2019/11/30/329348.html",
    "https://response.jp/article/2019/06/04/323073.html",
    "https://response.jp/article/2019/06/08/323236.html",
    "https://response.jp/article/2019/09/20/326717.html",
    "https://response.jp/article/2019/11/09/328553.html",
    "https://response.jp/article/2019/04/17/321474.html",
    "https://response.jp/article/2019/12/04/329454.html",
    "https://response.jp/article/2019/03/15/320170.html",
    "https://response.jp/article/2019/04/16/321415.html",
    "https://response.jp/article/2019/07/23/324728.html",
    "https://response.jp/special/recent/3532/%E4%BA%BA%E5%B7%A5%E7%9F%A5%E8%83%BD%EF%BC%88AI%EF%BC%89",
    "https://response.jp/article/2019/08/02/325155.html",
    "https://response.jp/article/2019/06/20/323636.html",
    "https://response.jp/article/2019/12/08/329583.html",
    "https://response.jp/article/2017/06/28/296736.html"
)

# Overwrite the existing data rows (A3:A20) with the refreshed URLs.
for ($i = 0; $i -lt $urls.Length; $i++) {
    $ws.Cells.Item($i + 3, 1).Value = $urls[$i]
}

# The previous data block ran through row 22; the refreshed block only
# needs through row 20, so drop the two now-unused trailing rows.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(21).Delete()

# Selection moves back to A3 after the refresh.
$ws.Range("A3").Select() | Out-Null
